$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: num_customers 47 -> 48 (cohort_size stays 2252), retention_rate recalculated
$ws.Range("C27").Value = 48
$ws.Range("E27").Value = 48 / 2252

# Row 37: num_customers and cohort_size 861 -> 865 (retention_rate stays 1)
$ws.Range("C37").Value = 865
$ws.Range("D37").Value = 865
